$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H40").Value = 6670762.5
$ws_ALC.Range("I40").Value = 16669941
$ws_ALC.Range("J40").Value = 4643.3335
$ws_ALC.Range("K40").Value = 16669941
$ws_ALC.Range("L40").Value = 4643.3335
$ws_ALC.Range("M40").Value = -16669766
$ws_ALC.Range("N40").Value = -4993.3335

$ws_ALC.Range("H51").Value = 127721200
$ws_ALC.Range("I51").Value = 340580860
$ws_ALC.Range("J51").Value = 5398.2
$ws_ALC.Range("K51").Value = 340580860
$ws_ALC.Range("L51").Value = 5398.2
$ws_ALC.Range("M51").Value = -340580376
$ws_ALC.Range("N51").Value = -6366.2

$ws_ALC.Range("H132").Value = 15876192
$ws_ALC.Range("I132").Value = 21742290
$ws_ALC.Range("J132").Value = 3217.1765
$ws_ALC.Range("K132").Value = 65226870
$ws_ALC.Range("L132").Value = 9651.529500000001
$ws_ALC.Range("M132").Value = -65224340
$ws_ALC.Range("N132").Value = -14711.5295

$ws_ALC.Range("H135").Value = 4427.8613
$ws_ALC.Range("I135").Value = 2690.4783
$ws_ALC.Range("J135").Value = 7501.6924
$ws_ALC.Range("K135").Value = 24214.3047
$ws_ALC.Range("L135").Value = 67515.2316
$ws_ALC.Range("M135").Value = -21679.3047
$ws_ALC.Range("N135").Value = -72585.2316

$ws_ALC.Range("H137").Value = 1795.5
$ws_ALC.Range("I137").Value = 1358.8334
$ws_ALC.Range("J137").Value = 2232.1667
$ws_ALC.Range("K137").Value = 4076.5002
$ws_ALC.Range("L137").Value = 6696.500100000001
$ws_ALC.Range("M137").Value = -1526.5002
$ws_ALC.Range("N137").Value = -11796.5001

$ws_ALC.Range("H138").Value = 353507.6
$ws_ALC.Range("I138").Value = 1500
$ws_ALC.Range("J138").Value = 394123.84
$ws_ALC.Range("K138").Value = 4500
$ws_ALC.Range("L138").Value = 1182371.52
$ws_ALC.Range("M138").Value = 640
$ws_ALC.Range("N138").Value = -1192651.52

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 2982.8484
$ws_ARM.Range("I32").Value = 2967
$ws_ARM.Range("K32").Value = 2967
$ws_ARM.Range("M32").Value = -2680

$ws_ARM.Range("H61").Value = 7518.091
$ws_ARM.Range("J61").Value = 2450
$ws_ARM.Range("L61").Value = 2450
$ws_ARM.Range("N61").Value = -2874

$ws_ARM.Range("H63").Value = 2348.5
$ws_ARM.Range("I63").Value = 2334.7273
$ws_ARM.Range("K63").Value = 2334.7273
$ws_ARM.Range("M63").Value = -1648.7273

$ws_ARM.Range("H66").Value = 2348.5
$ws_ARM.Range("I66").Value = 2334.7273
$ws_ARM.Range("K66").Value = 11673.6365
$ws_ARM.Range("M66").Value = -8241.636500000001

$ws_ARM.Range("H132").Value = 2986.4075
$ws_ARM.Range("I132").Value = 1749.9412
$ws_ARM.Range("J132").Value = 5088.4
$ws_ARM.Range("K132").Value = 5249.8236
$ws_ARM.Range("L132").Value = 15265.2
$ws_ARM.Range("M132").Value = -2719.8236
$ws_ARM.Range("N132").Value = -20325.2

$ws_ARM.Range("H136").Value = 7518.091
$ws_ARM.Range("J136").Value = 2450
$ws_ARM.Range("L136").Value = 7350
$ws_ARM.Range("N136").Value = -12450

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H82").Value = 10532.8125
$ws_BSM.Range("I82").Value = 2374.4285
$ws_BSM.Range("K82").Value = 2374.4285
$ws_BSM.Range("M82").Value = -1991.4285

$ws_BSM.Range("H85").Value = 10532.8125
$ws_BSM.Range("I85").Value = 2374.4285
$ws_BSM.Range("K85").Value = 2374.4285
$ws_BSM.Range("M85").Value = -1048.4285

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H16").Value = 2159.2307
$ws_CRP.Range("I16").Value = 2159.2307
$ws_CRP.Range("K16").Value = 2159.2307
$ws_CRP.Range("M16").Value = -1872.2307

$ws_CRP.Range("H22").Value = 527.75
$ws_CRP.Range("I22").Value = 496.33334
$ws_CRP.Range("K22").Value = 496.33334
$ws_CRP.Range("M22").Value = -146.33334

$ws_CRP.Range("H31").Value = 4506
$ws_CRP.Range("I31").Value = 2676.7
$ws_CRP.Range("J31").Value = 5044.0293
$ws_CRP.Range("K31").Value = 2676.7
$ws_CRP.Range("L31").Value = 5044.0293
$ws_CRP.Range("M31").Value = -2381.7
$ws_CRP.Range("N31").Value = -5634.0293

$ws_CRP.Range("H34").Value = 4506
$ws_CRP.Range("I34").Value = 2676.7
$ws_CRP.Range("J34").Value = 5044.0293
$ws_CRP.Range("K34").Value = 2676.7
$ws_CRP.Range("L34").Value = 5044.0293
$ws_CRP.Range("M34").Value = -2474.7
$ws_CRP.Range("N34").Value = -5448.0293

$ws_CRP.Range("H58").Value = 1883.6086
$ws_CRP.Range("I58").Value = 1943.4667
$ws_CRP.Range("J58").Value = 1771.375
$ws_CRP.Range("K58").Value = 1943.4667
$ws_CRP.Range("L58").Value = 1771.375
$ws_CRP.Range("M58").Value = -1740.4667
$ws_CRP.Range("N58").Value = -2177.375

$ws_CRP.Range("H113").Value = 2159.2307
$ws_CRP.Range("I113").Value = 2159.2307
$ws_CRP.Range("K113").Value = 2159.2307
$ws_CRP.Range("M113").Value = 10.76929999999993

$ws_CRP.Range("H132").Value = 1364.0571
$ws_CRP.Range("I132").Value = 1295.3226
$ws_CRP.Range("K132").Value = 3885.9678
$ws_CRP.Range("M132").Value = -1355.9678

$ws_CRP.Range("H136").Value = 1883.6086
$ws_CRP.Range("I136").Value = 1943.4667
$ws_CRP.Range("J136").Value = 1771.375
$ws_CRP.Range("K136").Value = 5830.4001
$ws_CRP.Range("L136").Value = 5314.125
$ws_CRP.Range("M136").Value = -3280.4001
$ws_CRP.Range("N136").Value = -10414.125

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H68").Value = 3835.7036
$ws_CUL.Range("I68").Value = 1665.3334
$ws_CUL.Range("J68").Value = 4920.8887
$ws_CUL.Range("K68").Value = 4996.0002
$ws_CUL.Range("L68").Value = 14762.6661
$ws_CUL.Range("M68").Value = -4185.0002
$ws_CUL.Range("N68").Value = -16384.6661

$ws_CUL.Range("H71").Value = 3835.7036
$ws_CUL.Range("I71").Value = 1665.3334
$ws_CUL.Range("J71").Value = 4920.8887
$ws_CUL.Range("K71").Value = 14988.0006
$ws_CUL.Range("L71").Value = 44287.99830000001
$ws_CUL.Range("M71").Value = -10932.0006
$ws_CUL.Range("N71").Value = -52399.99830000001

$ws_CUL.Range("H113").Value = 806.6177
$ws_CUL.Range("I113").Value = 779.3
$ws_CUL.Range("J113").Value = 818
$ws_CUL.Range("K113").Value = 2337.9
$ws_CUL.Range("L113").Value = 2454
$ws_CUL.Range("M113").Value = -167.8999999999996
$ws_CUL.Range("N113").Value = -6794

$ws_CUL.Range("H137").Value = 41673570
$ws_CUL.Range("I137").Value = 50004244
$ws_CUL.Range("J137").Value = 37045412
$ws_CUL.Range("K137").Value = 150012732
$ws_CUL.Range("L137").Value = 111136236
$ws_CUL.Range("M137").Value = -150007632
$ws_CUL.Range("N137").Value = -111146436

$ws_CUL.Range("H140").Value = 1764.4918
$ws_CUL.Range("I140").Value = 1246.9592
$ws_CUL.Range("J140").Value = 3877.75
$ws_CUL.Range("K140").Value = 3740.8776
$ws_CUL.Range("L140").Value = 11633.25
$ws_CUL.Range("M140").Value = 1439.1224
$ws_CUL.Range("N140").Value = -21993.25

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 3806.4666
$ws_LTW.Range("I7").Value = 3007
$ws_LTW.Range("J7").Value = 14999
$ws_LTW.Range("K7").Value = 3007
$ws_LTW.Range("L7").Value = 14999
$ws_LTW.Range("M7").Value = -2895
$ws_LTW.Range("N7").Value = -15223

$ws_LTW.Range("H126").Value = 3806.4666
$ws_LTW.Range("I126").Value = 3007
$ws_LTW.Range("J126").Value = 14999
$ws_LTW.Range("K126").Value = 9021
$ws_LTW.Range("L126").Value = 44997
$ws_LTW.Range("M126").Value = -6551
$ws_LTW.Range("N126").Value = -49937

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H96").Value = 4181.6665
$ws_WVR.Range("I96").Value = 3795
$ws_WVR.Range("J96").Value = 4375
$ws_WVR.Range("K96").Value = 3795
$ws_WVR.Range("L96").Value = 4375
$ws_WVR.Range("M96").Value = -2422
$ws_WVR.Range("N96").Value = -7121
